# Scheduled runner refresh: update Tonberry market-price columns (H/I/J/K/L/M/N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets with the latest
# Universalis price snapshot. Values only; no structural/formatting changes.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 76: Warding Off Temptation
$ws.Range("H76").Value = 3333
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()

# Row 79: The Garden of Arcane Delights (L)
$ws.Range("H79").Value = 3333
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()

# Row 92: Whinier than the Sword
$ws.Range("H92").Value = 14706067
$ws.Range("I92").Value = 15625134
$ws.Range("K92").Value = 15625134
$ws.Range("M92").Value = -15623886

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 819.06665
$ws.Range("I132").Value = 713.2174
$ws.Range("J132").Value = 1166.8572
$ws.Range("K132").Value = 2139.6522
$ws.Range("L132").Value = 3500.5716
$ws.Range("M132").Value = 390.3478
$ws.Range("N132").Value = -8560.571599999999

$ws = $wb.Worksheets.Item("ARM")
# Row 22: Kiss the Pan (Good-bye)
$ws.Range("H22").Value = 4999.923
$ws.Range("I22").Value = 4999.923
$ws.Range("K22").Value = 4999.923
$ws.Range("M22").Value = -4700.923

# Row 63: Rivets Run through It
$ws.Range("H63").Value = 3201
$ws.Range("I63").Value = 4201.3335
$ws.Range("J63").Value = 200
$ws.Range("K63").Value = 4201.3335
$ws.Range("L63").Value = 200
$ws.Range("M63").Value = -3515.3335
$ws.Range("N63").Value = -1572

# Row 66: A Riveting Revival (L)
$ws.Range("H66").Value = 3201
$ws.Range("I66").Value = 4201.3335
$ws.Range("J66").Value = 200
$ws.Range("K66").Value = 21006.6675
$ws.Range("L66").Value = 1000
$ws.Range("M66").Value = -17574.6675
$ws.Range("N66").Value = -7864

# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 2711.625
$ws.Range("I74").Value = 1000
$ws.Range("K74").Value = 1000
$ws.Range("M74").Value = -126

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 2711.625
$ws.Range("I77").Value = 1000
$ws.Range("K77").Value = 5000
$ws.Range("M77").Value = -632

# Row 88: The Mast Chance
$ws.Range("H88").Value = 4720
$ws.Range("I88").Value = 1800
$ws.Range("J88").Value = 6666.6665
$ws.Range("K88").Value = 1800
$ws.Range("L88").Value = 6666.6665
$ws.Range("M88").Value = -1394
$ws.Range("N88").Value = -7478.6665

# Row 91: The Rose and the Riveter (L)
$ws.Range("H91").Value = 4720
$ws.Range("I91").Value = 1800
$ws.Range("J91").Value = 6666.6665
$ws.Range("K91").Value = 1800
$ws.Range("L91").Value = 6666.6665
$ws.Range("M91").Value = -396
$ws.Range("N91").Value = -9474.666499999999

$ws = $wb.Worksheets.Item("BSM")
# Row 7: Thank You for Your Business
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

# Row 22: Riveting Run
$ws.Range("H22").Value = 747
$ws.Range("I22").Value = 747
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 747
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -574
$ws.Range("N22").ClearContents()

# Row 86: Through Thick and Thin
$ws.Range("H86").Value = 1718955.2
$ws.Range("I86").Value = 2006237.6
$ws.Range("K86").Value = 2006237.6
$ws.Range("M86").Value = -2005114.6

# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 1718955.2
$ws.Range("I89").Value = 2006237.6
$ws.Range("K89").Value = 10031188
$ws.Range("M89").Value = -10025572

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 4780.3667
$ws.Range("I134").Value = 5154.5386
$ws.Range("J134").Value = 2348.25
$ws.Range("K134").Value = 15463.6158
$ws.Range("L134").Value = 7044.75
$ws.Range("M134").Value = -12928.6158
$ws.Range("N134").Value = -12114.75

$ws = $wb.Worksheets.Item("CRP")
# Row 32: Daddy's Little Girl
$ws.Range("H32").Value = 14999.5
$ws.Range("I32").Value = 14999
$ws.Range("J32").Value = 15000
$ws.Range("K32").Value = 14999
$ws.Range("L32").Value = 15000
$ws.Range("M32").Value = -14683
$ws.Range("N32").Value = -15632

# Row 94: Beech, Please
$ws.Range("H94").Value = 1134.6364
$ws.Range("J94").Value = 1365.6666
$ws.Range("L94").Value = 1365.6666
$ws.Range("N94").Value = -2267.6666

# Row 105: Zelkova, My Love
$ws.Range("H105").Value = 1101.091
$ws.Range("I105").Value = 1126.8572
$ws.Range("J105").Value = 1056
$ws.Range("K105").Value = 1126.8572
$ws.Range("L105").Value = 1056
$ws.Range("M105").Value = 620.1428000000001
$ws.Range("N105").Value = -4550

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 2011.1333
$ws.Range("I132").Value = 1524.375
$ws.Range("J132").Value = 2567.4285
$ws.Range("K132").Value = 4573.125
$ws.Range("L132").Value = 7702.2855
$ws.Range("M132").Value = -2043.125
$ws.Range("N132").Value = -12762.2855

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap
$ws.Range("H5").Value = 427.57895
$ws.Range("I5").Value = 369.45456
$ws.Range("J5").Value = 507.5
$ws.Range("K5").Value = 1108.36368
$ws.Range("L5").Value = 1522.5
$ws.Range("M5").Value = -996.3636799999999
$ws.Range("N5").Value = -1746.5

# Row 12: Butter Me Up
$ws.Range("H12").Value = 117.35294
$ws.Range("I12").Value = 47.857143
$ws.Range("J12").Value = 166
$ws.Range("K12").Value = 143.571429
$ws.Range("L12").Value = 498
$ws.Range("M12").Value = 29.42857100000001
$ws.Range("N12").Value = -844

# Row 63: The Next to Last Supper
$ws.Range("H63").Value = 2676.6667

# Row 66: Nostalgia through the Stomach (L)
$ws.Range("H66").Value = 2676.6667

# Row 75: Breakfast of Champions
$ws.Range("H75").Value = 434.83334
$ws.Range("J75").Value = 467.8
$ws.Range("L75").Value = 1403.4
$ws.Range("N75").Value = -3399.4

# Row 78: Emerald Soup for the Soul (L)
$ws.Range("H78").Value = 434.83334
$ws.Range("J78").Value = 467.8
$ws.Range("L78").Value = 4210.2
$ws.Range("N78").Value = -14194.2

# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 49032.332
$ws.Range("J131").Value = 56458.08
$ws.Range("L131").Value = 169374.24
$ws.Range("N131").Value = -179454.24

# Row 135: Not-so-secret Ingredient
$ws.Range("H135").Value = 427.57895
$ws.Range("I135").Value = 369.45456
$ws.Range("J135").Value = 507.5
$ws.Range("K135").Value = 3325.09104
$ws.Range("L135").Value = 4567.5
$ws.Range("M135").Value = -790.0910400000002
$ws.Range("N135").Value = -9637.5

$ws = $wb.Worksheets.Item("GSM")
# Row 29: Music to Their Ears
$ws.Range("H29").Value = 72002
$ws.Range("I29").Value = 66664.664
$ws.Range("J29").Value = 80008
$ws.Range("K29").Value = 66664.664
$ws.Range("L29").Value = 80008
$ws.Range("M29").Value = -66374.664
$ws.Range("N29").Value = -80588

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 2174.9546
$ws.Range("I22").Value = 2225.9333
$ws.Range("J22").Value = 2065.7144
$ws.Range("K22").Value = 2225.9333
$ws.Range("L22").Value = 2065.7144
$ws.Range("M22").Value = -1930.9333
$ws.Range("N22").Value = -2655.7144

# Row 27: Fire and Hide
$ws.Range("H27").Value = 2174.9546
$ws.Range("I27").Value = 2225.9333
$ws.Range("J27").Value = 2065.7144
$ws.Range("K27").Value = 2225.9333
$ws.Range("L27").Value = 2065.7144
$ws.Range("M27").Value = -2118.9333
$ws.Range("N27").Value = -2279.7144

# Row 34: Breeches Served Cold
$ws.Range("H34").Value = 8420
$ws.Range("J34").Value = 8420
$ws.Range("L34").Value = 8420
$ws.Range("N34").Value = -8764

# Row 93: Hide to Go Seek
$ws.Range("H93").Value = 15152439
$ws.Range("I93").Value = 812.44446
$ws.Range("K93").Value = 812.44446
$ws.Range("M93").Value = 435.55554

# Row 134: Freezing Fingers
$ws.Range("H134").Value = 49428.75
$ws.Range("J134").Value = 49428.75
$ws.Range("L134").Value = 49428.75
$ws.Range("N134").Value = -59568.75

# Row 138: Freezing Toes
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# Row 141: Just Generally Freezing
$ws.Range("H141").Value = 45439.9
$ws.Range("J141").Value = 44933.89
$ws.Range("L141").Value = 44933.89
$ws.Range("N141").Value = -55293.89

$ws = $wb.Worksheets.Item("WVR")
# Row 62: Pride Up in Smoke
$ws.Range("H62").Value = 5108.1665
$ws.Range("I62").Value = 4862.25
$ws.Range("J62").Value = 5600
$ws.Range("K62").Value = 4862.25
$ws.Range("L62").Value = 5600
$ws.Range("M62").Value = -4238.25
$ws.Range("N62").Value = -6848

# Row 65: Desperate for Diversionaries (L)
$ws.Range("H65").Value = 5108.1665
$ws.Range("I65").Value = 4862.25
$ws.Range("J65").Value = 5600
$ws.Range("K65").Value = 24311.25
$ws.Range("L65").Value = 28000
$ws.Range("M65").Value = -21191.25
$ws.Range("N65").Value = -34240

# Row 133: Begin with the Basics
$ws.Range("H133").Value = 49999.5
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 49999.5
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 49999.5
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -60119.5
